# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) previously held a "Strike#"-style count; this
# recomputes/rewrites it with the new "K" values for every data row
# (rows 2-44) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,2,1,0,1,2,0,1,1,2,0,0,1,1,4,2,3,2,1,3,6,2,3,3,5,3,6,0,2,3,1,3,2,2,5,1,2,2,2,2,3,1,1)

$firstRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
